$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update the time_taken (F) column on the "data" sheet with refreshed query timestamps
$dataSheet.Range("F2").Value = "2021-10-05 14:21:13.274315"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:13.274324"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:13.274327"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:13.274330"
$dataSheet.Range("F6").Value = "2021-10-05 14:21:13.274333"
$dataSheet.Range("F7").Value = "2021-10-05 14:21:13.274336"
$dataSheet.Range("F8").Value = "2021-10-05 14:21:13.274339"
$dataSheet.Range("F9").Value = "2021-10-05 14:21:13.274341"
$dataSheet.Range("F10").Value = "2021-10-05 14:21:13.274344"
$dataSheet.Range("F11").Value = "2021-10-05 14:21:13.274347"
$dataSheet.Range("F12").Value = "2021-10-05 14:21:13.274349"
$dataSheet.Range("F13").Value = "2021-10-05 14:21:13.274352"
$dataSheet.Range("F14").Value = "2021-10-05 14:21:13.274355"
$dataSheet.Range("F15").Value = "2021-10-05 14:21:13.274358"
$dataSheet.Range("F16").Value = "2021-10-05 14:21:13.274360"
$dataSheet.Range("F17").Value = "2021-10-05 14:21:13.274363"
$dataSheet.Range("F18").Value = "2021-10-05 14:21:13.274366"
$dataSheet.Range("F19").Value = "2021-10-05 14:21:13.274369"
$dataSheet.Range("F20").Value = "2021-10-05 14:21:13.274371"
$dataSheet.Range("F21").Value = "2021-10-05 14:21:13.274374"
$dataSheet.Range("F22").Value = "2021-10-05 14:21:13.274377"
$dataSheet.Range("F23").Value = "2021-10-05 14:21:13.274380"
$dataSheet.Range("F24").Value = "2021-10-05 14:21:13.274382"
$dataSheet.Range("F25").Value = "2021-10-05 14:21:13.274385"
$dataSheet.Range("F26").Value = "2021-10-05 14:21:13.274388"
$dataSheet.Range("F27").Value = "2021-10-05 14:21:13.274391"

# Add the new "metadata" sheet right after "data"
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row
$ws.Range("A2").Value = 0
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160

$ws.Range("B2").Value = "Insulin resistance (including lipodystrophy)"
$ws.Range("C2").Value = 174
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.13"
$ws.Range("E2").Value = "2021-07-28T09:57:34.668077Z"
$ws.Range("F2").Value = "2021-10-05 14:21:13.270668"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/174/?format=json"

# Keep "data" as the active/selected sheet (unchanged bookViews semantics)
$dataSheet.Activate()

Write-Host "metadata sheet added; time_taken column refreshed"
